$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph reads "ตาราง … Sequence Diagram" (Table ... Sequence Diagram).
# It must become "ตารางที่ 1 Sequence Diagram" (Table No. 1 Sequence Diagram):
#   * the run holding the lone space right after "ตาราง" becomes "ที่ "
#   * the run holding the ellipsis "… " is replaced by a run with just "1"
#     followed by a brand-new run holding a single space " "
#   * "Sequence" / " Diagram" stay as-is
# ------------------------------------------------------------------

# Anchor on the literal "ตาราง" text so we do not depend on hard-coded offsets.
$anchor = $d.Content
$anchor.Find.Execute("ตาราง", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$tabEnd = $anchor.End

# --- 1) turn the single space run after "ตาราง" into "ที่ " ---
$spaceRun = $d.Range($tabEnd, $tabEnd + 1)
$spaceRun.Text = "ที่ "

# Re-split the new text from "ตาราง" so it keeps living in its own run
# (a plain property round-trip is enough to stop the engine re-merging
# it with the identically formatted run to its left).
$newSpaceRun = $d.Range($tabEnd, $tabEnd + 4)
$newSpaceRun.Font.Bold = $true
$newSpaceRun.Font.Bold = $false

# --- 2) turn the ellipsis run "… " into "1" + a fresh " " run ---
$ellipsisStart = $tabEnd + 4
$ellipsisRun = $d.Range($ellipsisStart, $ellipsisStart + 2)
$ellipsisRun.Text = "1"

# Re-split "1" away from the text on its left.
$oneRun = $d.Range($ellipsisStart, $ellipsisStart + 1)
$oneRun.Font.Bold = $true
$oneRun.Font.Bold = $false

# Insert the brand-new standalone space run right after "1".
$insertPoint = $d.Range($ellipsisStart + 1, $ellipsisStart + 1)
$insertPoint.InsertBefore(" ")

# Re-split the freshly inserted space from "1".
$newSpace2 = $d.Range($ellipsisStart + 1, $ellipsisStart + 2)
$newSpace2.Font.Bold = $true
$newSpace2.Font.Bold = $false

# Re-split "Sequence" away from the new space run so it (and " Diagram"
# behind it) keep standing on their own, unaffected by the edit above.
$seqRun = $d.Range($ellipsisStart + 2, $ellipsisStart + 10)
$seqRun.Font.Bold = $true
$seqRun.Font.Bold = $false
